# onCampList shows all staff on camp now
# Adds a new "04-29-2022" sheet (a filled-in copy of the "Daily Attendance
# Template" sheet) between "Key" and "Daily Attendance Template", and marks
# one staff member as having left camp on the Key sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new dated sheet by duplicating the template, right after
#    "Key" (matches sheetId ordering: Key, 04-29-2022, Daily Attendance
#    Template).
# ---------------------------------------------------------------------
$keySheet = $wb.Worksheets.Item("Key")
$templateSheet = $wb.Worksheets.Item("Daily Attendance Template")
$templateSheet.Copy($null, $keySheet)

$newSheet = $wb.Worksheets.Item("Daily Attendance Template (2)")
$newSheet.Name = "04-29-2022"

# The template sheet is protected, so unlock it before editing.
$newSheet.Unprotect()

# ---------------------------------------------------------------------
# 2. Row 2: Staff Member 4 (Bunk 2) leaving camp at 4:27 PM.
# ---------------------------------------------------------------------
$newSheet.Range("A2").Value2 = "Bunk 2"
$newSheet.Range("B2").Value2 = "Staff Member 4"

$newSheet.Range("C2").Value2 = "Staff Member 4 ID"
$newSheet.Range("C2").Borders.Item(10).LineStyle = 1

$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("D2").Value2 = "4:27 PM"

$newSheet.Range("E2").Value2 = "Leaving`nCamp"
$newSheet.Range("E2").Interior.Color = 10066410

$newSheet.Range("I2").Value2 = "1:00 AM"

# ---------------------------------------------------------------------
# 3. Row 3: Staff Member 4 (Bunk 2) returning to camp at 4:30 PM.
# ---------------------------------------------------------------------
$newSheet.Range("A3").Value2 = "Bunk 2"
$newSheet.Range("B3").Value2 = "Staff Member 4"

$newSheet.Range("C3").Value2 = "Staff Member 4 ID"
$newSheet.Range("C3").Borders.Item(10).LineStyle = 1

$newSheet.Range("D3").Style = "Normal"
$newSheet.Range("D3").Value2 = "4:29 PM"

$newSheet.Range("E3").Value2 = "4:30 PM"
$newSheet.Range("E3").Interior.Color = 13492663

$newSheet.Range("I3").Value2 = "1:00 AM"

# ---------------------------------------------------------------------
# 4. Curfew / summary panel on the right (G:I).
# ---------------------------------------------------------------------
$newSheet.Range("I4").Value2 = "5:00 PM"

$newSheet.Range("I6").Value2 = 1
$newSheet.Range("I7").Value2 = 1
$newSheet.Range("I8").Value2 = 0
$newSheet.Range("I10").Value2 = 0

# Re-lock the sheet the same way the template was protected.
$newSheet.Protect()

# ---------------------------------------------------------------------
# 5. Key sheet: mark Staff Member 4 ID row with the "on camp list" flag.
# ---------------------------------------------------------------------
$keySheet.Range("D5").Value2 = 1
